# Generate Report for Handback
# Replace the first handback file's identifiers/timestamps with the new run's,
# and likewise for the second file. Update hyperlink display text to match.

$wb = $excel.ActiveWorkbook

$oldFile1 = "bc481864-bdb7-4409-a64b-0dbffbf8778a.md"
$newFile1 = "6f6002af-4bac-4223-b75b-3cc77185eb73.md"

$oldFile2 = "d4073668-b61f-49ab-9751-4aa41cdea716.md"
$newFile2 = "ffffb6fe7995-e509-4c66-87ad-b74a8e7c687f.md"

$newPath1 = "e2e\$newFile1"
$newPath2 = "e2e\$newFile2"

$newXlf1zhcn = "6f6002af-4bac-4223-b75b-3cc77185eb73.53e624c6b7227c39bc612a23d91d7edc86c7f095.zh-cn.xlf"
$newXlf1dede = "6f6002af-4bac-4223-b75b-3cc77185eb73.53e624c6b7227c39bc612a23d91d7edc86c7f095.de-de.xlf"

$newOverviewDate = "2016-08-30 19:17:37"

$newZhCnHandoffDate = "2016-08-30 19:17:31"
$newZhCnHandbackDate = "2016-08-30 19:17:54"

$newDeDeHandbackDate = "2016-08-30 19:18:03"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newPath1
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newPath2
$wsOverview.Range("G3").Value = $newOverviewDate

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = $newPath1
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = $newPath2
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("G2").Value = $newXlf1zhcn
$wsZhCn.Range("H2").Value = $newZhCnHandoffDate
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("J2").Value = $newXlf1zhcn
$wsZhCn.Range("K2").Value = $newZhCnHandbackDate

$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("G3").Value = $newXlf1zhcn
$wsZhCn.Range("H3").Value = $newZhCnHandoffDate
$wsZhCn.Range("I3").Value = $newFile2
$wsZhCn.Range("J3").Value = $newXlf1zhcn
$wsZhCn.Range("K3").Value = $newZhCnHandbackDate

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFile1
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = $newFile1
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $newFile2
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = $newFile2
    }
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("G2").Value = $newXlf1dede
$wsDeDe.Range("H2").Value = $newOverviewDate
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("J2").Value = $newXlf1dede
$wsDeDe.Range("K2").Value = $newDeDeHandbackDate

$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("G3").Value = $newXlf1dede
$wsDeDe.Range("H3").Value = $newOverviewDate
$wsDeDe.Range("I3").Value = $newFile2
$wsDeDe.Range("J3").Value = $newXlf1dede
$wsDeDe.Range("K3").Value = $newDeDeHandbackDate

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFile1
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = $newFile1
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $newFile2
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = $newFile2
    }
}
